$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of survey data (row 4)
$ws.Range("A4").Value = "Analysis of Inter-Satellite Link Paths for LEO Mega-Constellation Networks"
$ws.Range("B4").Value = "TVT"
$ws.Range("C4").Value = "1.提出通过经纬度计算估算两个卫星之间的ISL跳数，其实可以通过给卫星编号来计算。。。。`n2."

# Match the formatting used by the existing rows: column C wraps text,
# vertically centered (same style as C3/D3/E3), and the row is tall enough
# to show the wrapped paragraph.
$ws.Range("C4").WrapText = $true
$ws.Range("C4").VerticalAlignment = -4108
$ws.Rows.Item(4).RowHeight = 112.5

# Move the active selection, matching where the author left off editing.
$ws.Range("C10").Select()
